$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 47, shifting existing rows 47-109 down to 49-111
$ws.Rows.Item(47).EntireRow.Insert()
$ws.Rows.Item(47).EntireRow.Insert()

# New row 47: Terminal La Palmera de La Serena - Damasco - Dina
$ws.Cells.Item(47,1).Value = 8
$ws.Cells.Item(47,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47,3).Value = "Coquimbo"
$ws.Cells.Item(47,4).Value = 44935
$ws.Cells.Item(47,5).Value = 4
$ws.Cells.Item(47,6).Value = "Fruta"
$ws.Cells.Item(47,7).Value = 100103
$ws.Cells.Item(47,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(47,9).Value = 100103003
$ws.Cells.Item(47,10).Value = "Damasco"
$ws.Cells.Item(47,11).Value = "Dina"
$ws.Cells.Item(47,12).Value = "Primera"
$ws.Cells.Item(47,13).Value = 200
$ws.Cells.Item(47,14).Value = 19000
$ws.Cells.Item(47,15).Value = 20000
$ws.Cells.Item(47,16).Value = 19500
$ws.Cells.Item(47,17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(47,18).Value = "Región Metropolitana"
$ws.Cells.Item(47,19).Value = 1219
$ws.Cells.Item(47,20).Value = 16

# New row 48: Terminal La Palmera de La Serena - Damasco - Dina
$ws.Cells.Item(48,1).Value = 8
$ws.Cells.Item(48,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(48,3).Value = "Coquimbo"
$ws.Cells.Item(48,4).Value = 44935
$ws.Cells.Item(48,5).Value = 4
$ws.Cells.Item(48,6).Value = "Fruta"
$ws.Cells.Item(48,7).Value = 100103
$ws.Cells.Item(48,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(48,9).Value = 100103003
$ws.Cells.Item(48,10).Value = "Damasco"
$ws.Cells.Item(48,11).Value = "Dina"
$ws.Cells.Item(48,12).Value = "Segunda"
$ws.Cells.Item(48,13).Value = 160
$ws.Cells.Item(48,14).Value = 15000
$ws.Cells.Item(48,15).Value = 16000
$ws.Cells.Item(48,16).Value = 15500
$ws.Cells.Item(48,17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(48,18).Value = "Región Metropolitana"
$ws.Cells.Item(48,19).Value = 969
$ws.Cells.Item(48,20).Value = 16
